# Applies per-cell value updates produced by the scheduled Chocobo_Profits
# market-price refresh (H:N = price/profit columns on each crafting-job sheet).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2763
$ws.Range("I98").Value = 893.05554
$ws.Range("J98").Value = 7571.4287
$ws.Range("K98").Value = 893.05554
$ws.Range("L98").Value = 7571.4287
$ws.Range("M98").Value = 604.94446
$ws.Range("N98").Value = -10567.4287
$ws.Range("H100").Value = 12501744
$ws.Range("I100").Value = 12501744
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 12501744
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -12501203
$ws.Range("N100").ClearContents()
$ws.Range("H116").Value = 339896.44
$ws.Range("I116").Value = 591240.75
$ws.Range("J116").Value = 11215.385
$ws.Range("K116").Value = 591240.75
$ws.Range("L116").Value = 11215.385
$ws.Range("M116").Value = -587798.75
$ws.Range("N116").Value = -18099.385
$ws.Range("H122").Value = 2763
$ws.Range("I122").Value = 893.05554
$ws.Range("J122").Value = 7571.4287
$ws.Range("K122").Value = 2679.16662
$ws.Range("L122").Value = 22714.2861
$ws.Range("M122").Value = -229.16662
$ws.Range("N122").Value = -27614.2861
$ws.Range("H133").Value = 46823.637
$ws.Range("J133").Value = 46823.637
$ws.Range("L133").Value = 46823.637
$ws.Range("N133").Value = -56943.637
$ws.Range("H137").Value = 3465.0435
$ws.Range("I137").Value = 1420.4
$ws.Range("J137").Value = 7298.75
$ws.Range("K137").Value = 4261.200000000001
$ws.Range("L137").Value = 21896.25
$ws.Range("M137").Value = -1711.200000000001
$ws.Range("N137").Value = -26996.25
$ws.Range("H138").Value = 1898.62
$ws.Range("I138").Value = 672.5143
$ws.Range("J138").Value = 2558.8308
$ws.Range("K138").Value = 2017.5429
$ws.Range("L138").Value = 7676.492400000001
$ws.Range("M138").Value = 3122.4571
$ws.Range("N138").Value = -17956.4924

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3822.8276
$ws.Range("I74").Value = 3512.889
$ws.Range("J74").Value = 8007
$ws.Range("K74").Value = 3512.889
$ws.Range("L74").Value = 8007
$ws.Range("M74").Value = -2638.889
$ws.Range("N74").Value = -9755
$ws.Range("H77").Value = 3822.8276
$ws.Range("I77").Value = 3512.889
$ws.Range("J77").Value = 8007
$ws.Range("K77").Value = 17564.445
$ws.Range("L77").Value = 40035
$ws.Range("M77").Value = -13196.445
$ws.Range("N77").Value = -48771
$ws.Range("H132").Value = 2124.4
$ws.Range("I132").Value = 1060.0358
$ws.Range("J132").Value = 4607.9165
$ws.Range("K132").Value = 3180.1074
$ws.Range("L132").Value = 13823.7495
$ws.Range("M132").Value = -650.1074000000003
$ws.Range("N132").Value = -18883.7495

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1620.76
$ws.Range("I105").Value = 1622.9131
$ws.Range("J105").Value = 1596
$ws.Range("K105").Value = 1622.9131
$ws.Range("L105").Value = 1596
$ws.Range("M105").Value = 124.0869
$ws.Range("N105").Value = -5090
$ws.Range("H107").Value = 1999.3334
$ws.Range("I107").Value = 1999.3334
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1999.3334
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -79.33339999999998
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 1582.3334
$ws.Range("I134").Value = 973.9394
$ws.Range("J134").Value = 3255.4167
$ws.Range("K134").Value = 2921.8182
$ws.Range("L134").Value = 9766.250100000001
$ws.Range("M134").Value = -386.8181999999997
$ws.Range("N134").Value = -14836.2501

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1749.2179
$ws.Range("I58").Value = 1542.4615
$ws.Range("J58").Value = 2783
$ws.Range("K58").Value = 1542.4615
$ws.Range("L58").Value = 2783
$ws.Range("M58").Value = -1339.4615
$ws.Range("N58").Value = -3189
$ws.Range("H122").Value = 2862.6843
$ws.Range("I122").Value = 2188.5833
$ws.Range("J122").Value = 4018.2856
$ws.Range("K122").Value = 6565.749899999999
$ws.Range("L122").Value = 12054.8568
$ws.Range("M122").Value = -4115.749899999999
$ws.Range("N122").Value = -16954.8568
$ws.Range("H132").Value = 3363.5898
$ws.Range("I132").Value = 2945.5
$ws.Range("J132").Value = 4199.769
$ws.Range("K132").Value = 8836.5
$ws.Range("L132").Value = 12599.307
$ws.Range("M132").Value = -6306.5
$ws.Range("N132").Value = -17659.307
$ws.Range("H134").Value = 4169.6113
$ws.Range("I134").Value = 5277.909
$ws.Range("J134").Value = 2428
$ws.Range("K134").Value = 15833.727
$ws.Range("L134").Value = 7284
$ws.Range("M134").Value = -13298.727
$ws.Range("N134").Value = -12354
$ws.Range("H136").Value = 1749.2179
$ws.Range("I136").Value = 1542.4615
$ws.Range("J136").Value = 2783
$ws.Range("K136").Value = 4627.3845
$ws.Range("L136").Value = 8349
$ws.Range("M136").Value = -2077.3845
$ws.Range("N136").Value = -13449

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 742.61536
$ws.Range("J113").Value = 887
$ws.Range("L113").Value = 2661
$ws.Range("N113").Value = -7001
$ws.Range("H131").Value = 5814817
$ws.Range("J131").Value = 906.575
$ws.Range("L131").Value = 2719.725
$ws.Range("N131").Value = -12799.725

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3224.48
$ws.Range("I126").Value = 2992
$ws.Range("J126").Value = 4283.5557
$ws.Range("K126").Value = 8976
$ws.Range("L126").Value = 12850.6671
$ws.Range("M126").Value = -6506
$ws.Range("N126").Value = -17790.6671
$ws.Range("H132").Value = 2503
$ws.Range("I132").Value = 1365.5217
$ws.Range("J132").Value = 4683.1665
$ws.Range("K132").Value = 4096.5651
$ws.Range("L132").Value = 14049.4995
$ws.Range("M132").Value = -1566.5651
$ws.Range("N132").Value = -19109.4995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 309.2
$ws.Range("I55").Value = 299
$ws.Range("J55").Value = 350
$ws.Range("K55").Value = 299
$ws.Range("L55").Value = 350
$ws.Range("M55").Value = -126
$ws.Range("N55").Value = -696
$ws.Range("H100").Value = 1621.9166
$ws.Range("I100").Value = 1423.2778
$ws.Range("K100").Value = 1423.2778
$ws.Range("M100").Value = -882.2778000000001
$ws.Range("H122").Value = 4478.6665
$ws.Range("I122").Value = 2075
$ws.Range("J122").Value = 7225.7144
$ws.Range("K122").Value = 6225
$ws.Range("L122").Value = 21677.1432
$ws.Range("M122").Value = -3775
$ws.Range("N122").Value = -26577.1432

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6290625
$ws.Range("I132").Value = 808.4
$ws.Range("J132").Value = 41670844
$ws.Range("K132").Value = 2425.2
$ws.Range("L132").Value = 125012532
$ws.Range("M132").Value = 104.8000000000002
$ws.Range("N132").Value = -125017592

